# Fill in the home/away scores for the games that finished after the
# workbook was first uploaded (rows 14, 15, 16, 17 of Sheet1).
#
# Columns: A=game_id, B=home_team, C=away_team, D=home_score, E=away_score

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("D14").Value = 17
$ws.Range("E14").Value = 20

# Row 15
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 22

# Row 16 already has a home_score (D16); only away_score (E16) is new
$ws.Range("E16").Value = 20

# Row 17
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 20

# Match the numeric-score formatting already used on this sheet (e.g. D16)
# rather than leaving the new cells on the sheet's default/text style.
$ws.Range("D14:E15").Font.Name = "Calibri"
$ws.Range("E16").Font.Name = "Calibri"
$ws.Range("D17:E17").Font.Name = "Calibri"
